$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 123, shifting existing rows 123:152 down to 124:153
$ws.Rows("123:123").Insert()

# Populate the newly inserted row 123 with the new weekly record
$ws.Range("A123").Value = 7
$ws.Range("B123").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C123").Value = "Ñuble"
$ws.Range("D123").Value = 44943
$ws.Range("D123").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E123").Value = 16
$ws.Range("F123").Value = 100112040
$ws.Range("G123").Value = "Cilantro"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 300
$ws.Range("K123").Value = 700
$ws.Range("L123").Value = 800
$ws.Range("M123").Value = 750
$ws.Range("N123").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O123").Value = "Provincia de Diguillín"
$ws.Range("P123").Value = 750
$ws.Range("Q123").Value = 1
$ws.Range("R123").Value = "Hortaliza"
